# Fruta / hortaliza, semanal
# Insert a new data row at row 94 (shifting the existing rows 94:189 down to
# 95:190) and populate it with the new weekly price-reporting data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 94..189 down to 95..190, duplicating formatting of row 94.
$ws.Rows.Item(94).Insert()

$newRow = 94

$ws.Cells.Item($newRow, 1).Value2 = 7
$ws.Cells.Item($newRow, 2).Value2 = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item($newRow, 3).Value2 = 'Ñuble'
$ws.Cells.Item($newRow, 4).Value2 = 44539
$ws.Cells.Item($newRow, 5).Value2 = 16
$ws.Cells.Item($newRow, 6).Value2 = 100112009
$ws.Cells.Item($newRow, 7).Value2 = 'Acelga'
$ws.Cells.Item($newRow, 8).Value2 = 'Sin especificar'
$ws.Cells.Item($newRow, 9).Value2 = 'Primera'
$ws.Cells.Item($newRow, 10).Value2 = 120
$ws.Cells.Item($newRow, 11).Value2 = 350
$ws.Cells.Item($newRow, 12).Value2 = 400
$ws.Cells.Item($newRow, 13).Value2 = 375
$ws.Cells.Item($newRow, 14).Value2 = '$/atado 0,5 a 1 kilo'
$ws.Cells.Item($newRow, 15).Value2 = 'Provincia de Diguillín'
$ws.Cells.Item($newRow, 16).Value2 = 375
$ws.Cells.Item($newRow, 17).Value2 = 1
$ws.Cells.Item($newRow, 18).Value2 = 'Hortaliza'
